$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.315551519393921
$ws.Range("B1").Value = 1.086208462715149
$ws.Range("C1").Value = 4.003510475158691
$ws.Range("D1").Value = 4.320512771606445
$ws.Range("E1").Value = 0.8162689208984375
